$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: restyle existing row (s=4/5 -> s=6/7) and add empty A33 ---
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("C30:E30").Copy() | Out-Null
$ws.Range("C33:E33").PasteSpecial(-4122) | Out-Null

# --- Row 34: new row (style like row 29, ht 43.2) ---
$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A34:E34").PasteSpecial(-4122) | Out-Null

# --- Row 35: new row (style like row 30, default height) ---
$ws.Range("A30:E30").Copy() | Out-Null
$ws.Range("A35:E35").PasteSpecial(-4122) | Out-Null

# --- Row 36: new row (style like row 29, ht 43.2) ---
$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A36:E36").PasteSpecial(-4122) | Out-Null

# --- Row 37: new row (style like row 29, ht 43.2) ---
$ws.Range("A29:E29").Copy() | Out-Null
$ws.Range("A37:E37").PasteSpecial(-4122) | Out-Null

# --- Row 38: new row (style like row 32, ht 21.6, no A cell) ---
$ws.Range("B32:E32").Copy() | Out-Null
$ws.Range("B38:E38").PasteSpecial(-4122) | Out-Null

# --- Row 39: new row (style like row 18, default height, no A cell) ---
$ws.Range("B18:E18").Copy() | Out-Null
$ws.Range("B39:E39").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Set values, in the exact order new shared strings are introduced ---
$ws.Range("C34").Value = ' Oh my gosh! I got to meet the\nfamous [CS:N]Dusknoir[CR]!'
$ws.Range("C35").Value = ' It\''s such an honor! Yippee!'
$ws.Range("A34").Value = 'SCRIPT/T01P01A/um1101.ssb'
$ws.Range("D34").Value = ' О боже мой! Я встретилась с\nизвестным [CS:N]Даскнуаром[CR]!'
$ws.Range("D35").Value = ' Это такая честь! Ураааа!'
$ws.Range("E34").Value = ' Ï áïçå íïê! Ÿ âòóñåóéìàòû ò\néèâåòóîúí [CS:N]Äàòëîôàñïí[CR]!'
$ws.Range("E35").Value = ' Üóï óàëàÿ œåòóû! Ôñàààà!'
$ws.Range("A36").Value = 'SCRIPT/G01P03A/um1108.ssb'
$ws.Range("C36").Value = ' I had the chance to chat with\nthe great [CS:N]Dusknoir[CR] the other day.'
$ws.Range("C37").Value = ' I must say that [CS:N]Dusknoir[CR] is\nincredibly worldly and wise.'
$ws.Range("C38").Value = ' He\''s a wealth of hints and tips\nfor exploring too.'
$ws.Range("C39").Value = ' He\''s so very inspiring!'
$ws.Range("D36").Value = ' Недавно я смогла пообщаться\nс [CS:N]Даскнуаром[CR].'
$ws.Range("D37").Value = ' Должна сказать, что [CS:N]Даскнуар[CR]\nневероятно открытый и мудрый.'
$ws.Range("D38").Value = ' Ещё он знает множество моментов,\nполезных для исследований.'
$ws.Range("D39").Value = ' Он так меня вдохновляет!'
$ws.Range("E36").Value = ' Îåäàâîï ÿ òíïãìà ðïïáþàóûòÿ\nò [CS:N]Äàòëîôàñïí[CR].'
$ws.Range("E37").Value = ' Äïìçîà òëàèàóû, œóï [CS:N]Äàòëîôàñ[CR]\nîåâåñïÿóîï ïóëñúóúê é íôäñúê.'
$ws.Range("E38").Value = ' Åþæ ïî èîàåó íîïçåòóâï íïíåîóïâ,\nðïìåèîúö äìÿ éòòìåäïâàîéê.'
$ws.Range("E39").Value = ' Ïî óàë íåîÿ âäïöîïâìÿåó!'
$ws.Range("A37").Value = 'SCRIPT/G01P03A/um1113.ssb'

# --- Set line-number (column B) values ---
$ws.Range("B34").Value = 533
$ws.Range("B35").Value = 536
$ws.Range("B36").Value = 505
$ws.Range("B37").Value = 508
$ws.Range("B38").Value = 511
$ws.Range("B39").Value = 514

# --- Row heights for new rows ---
$ws.Rows.Item(34).RowHeight = 43.2
$ws.Rows.Item(36).RowHeight = 43.2
$ws.Rows.Item(37).RowHeight = 43.2
$ws.Rows.Item(38).RowHeight = 21.6

# --- Selection / active cell to match final view ---
$ws.Range("C39").Select() | Out-Null
